$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999764042957,
    0.99840872636903522,
    0.99444185349559966,
    0.99737906000968635,
    0.98716331046137751,
    0.96386655402163179,
    0.96154300311093077,
    0.95795443442491424,
    0.95616691640727103,
    0.94561898153384183,
    0.94405838322734748,
    0.94127889080847793,
    0.92999291149349372,
    0.92582571425501281,
    0.92323428172701205,
    0.92072783560568472,
    0.91701992044097724,
    0.91591098749826327,
    0.99720778719193892,
    0.98857517913633464,
    0.98717665980794256,
    0.98591215011157096,
    0.98740766877049668,
    0.97438743041662068,
    0.96793053037325727,
    0.96064648256160368,
    0.95640541301109039,
    0.93864529777759798,
    0.92610509363614391,
    0.92054821675901588,
    0.9128949141895899,
    0.91121556409244442,
    0.91069554208667858
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
